$wb = $excel.ActiveWorkbook

$wsConstants = $wb.Worksheets.Item("constants")
$wsTimeVariants = $wb.Worksheets.Item("time_variants")

# Delete the econ_*_shortcourse_mdr rows (28-32) first so row numbers above
# are unaffected while we do so.
$wsConstants.Rows("28:32").Delete() | Out-Null

# Delete the age_breakpoints row (row 6)
$wsConstants.Rows("6").Delete() | Out-Null

# Update selections / active sheet to roughly match the target workbook state
$wsConstants.Range("E25").Select() | Out-Null
$wsTimeVariants.Range("AR15").Select() | Out-Null

$wsConstants.Activate() | Out-Null

$wb.Save()
